# OP_148_SoliGrupEdad.xlsx - "added new data for Adm. de Tribunales"
#
# 1) Sheet "2020-2021": fix a mis-pointed header string in C1 (it had a stray
#    leading-space variant of "20 años o menos") and remove the three rows
#    that don't belong on this tab (Ponce / San Juan / Utuado), which also
#    shifts the Total row up and re-bases its SUM() formulas.
# 2) Add a new sheet "2024-2025" (after "2023-2024") with a full year's worth
#    of data, and make it the active tab/selection like the source file.

$wb = $excel.ActiveWorkbook

# --- 1) Fix up "2020-2021" -------------------------------------------------
$ws1 = $wb.Worksheets.Item("2020-2021")

# Correct header text (shared-string clean up happens automatically on save).
$ws1.Range("C1").Value = "20 años o menos"

# Remove Ponce / San Juan / Utuado rows (rows 12-14); Total (row 15) shifts up
# to row 12 and its SUM formulas automatically re-base to B2:B11 etc.
$ws1.Rows("12:14").Delete()

# --- 2) Add new sheet "2024-2025" with this year's data --------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "2024-2025"

$headers = @("Región","Cantidad total","20 años o menos","21-29 años","30-39 años","40-49 años","50-59 años","60 años o más","No indica")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$ws.Cells.Item(2, 1).Value = "Aguadilla"
$ws.Cells.Item(2, 2).Value = 7
$ws.Cells.Item(2, 3).Value = 6
$ws.Cells.Item(2, 8).Value = 1

$ws.Cells.Item(3, 1).Value = "Aibonito"
$ws.Cells.Item(3, 2).Value = 3
$ws.Cells.Item(3, 3).Value = 1
$ws.Cells.Item(3, 5).Value = 2

$ws.Cells.Item(4, 1).Value = "Arecibo"
$ws.Cells.Item(4, 2).Value = 14
$ws.Cells.Item(4, 3).Value = 6
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 9).Value = 5

$ws.Cells.Item(5, 1).Value = "Bayamón"
$ws.Cells.Item(5, 2).Value = 18
$ws.Cells.Item(5, 3).Value = 7
$ws.Cells.Item(5, 4).Value = 3
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 9).Value = 7

$ws.Cells.Item(6, 1).Value = "Caguas"
$ws.Cells.Item(6, 2).Value = 19
$ws.Cells.Item(6, 3).Value = 3
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 9).Value = 14

$ws.Cells.Item(7, 1).Value = "Carolina"
$ws.Cells.Item(7, 2).Value = 5
$ws.Cells.Item(7, 3).Value = 2
$ws.Cells.Item(7, 4).Value = 2
$ws.Cells.Item(7, 9).Value = 1

$ws.Cells.Item(8, 1).Value = "Fajardo"
$ws.Cells.Item(8, 2).Value = 1
$ws.Cells.Item(8, 6).Value = 1

$ws.Cells.Item(9, 1).Value = "Guayama"
$ws.Cells.Item(9, 2).Value = 3
$ws.Cells.Item(9, 4).Value = 1
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 1

$ws.Cells.Item(10, 1).Value = "Humacao"
$ws.Cells.Item(10, 2).Value = 4
$ws.Cells.Item(10, 4).Value = 2
$ws.Cells.Item(10, 6).Value = 2

$ws.Cells.Item(11, 1).Value = "Mayagüez"
$ws.Cells.Item(11, 2).Value = 4
$ws.Cells.Item(11, 3).Value = 3
$ws.Cells.Item(11, 6).Value = 1

$ws.Cells.Item(12, 1).Value = "Ponce"
$ws.Cells.Item(12, 2).Value = 19
$ws.Cells.Item(12, 3).Value = 13
$ws.Cells.Item(12, 4).Value = 1
$ws.Cells.Item(12, 5).Value = 2
$ws.Cells.Item(12, 8).Value = 1
$ws.Cells.Item(12, 9).Value = 2

$ws.Cells.Item(13, 1).Value = "San Juan"
$ws.Cells.Item(13, 2).Value = 5
$ws.Cells.Item(13, 4).Value = 1
$ws.Cells.Item(13, 5).Value = 2
$ws.Cells.Item(13, 9).Value = 2

$ws.Cells.Item(14, 1).Value = "Utuado"
$ws.Cells.Item(14, 2).Value = 7
$ws.Cells.Item(14, 3).Value = 3
$ws.Cells.Item(14, 4).Value = 1
$ws.Cells.Item(14, 6).Value = 2
$ws.Cells.Item(14, 7).Value = 1

$ws.Cells.Item(15, 1).Value = "Total"
$ws.Cells.Item(15, 2).Value = 109
$ws.Cells.Item(15, 3).Value = 44
$ws.Cells.Item(15, 4).Value = 11
$ws.Cells.Item(15, 5).Value = 11
$ws.Cells.Item(15, 6).Value = 8
$ws.Cells.Item(15, 7).Value = 2
$ws.Cells.Item(15, 8).Value = 2
$ws.Cells.Item(15, 9).Value = 31

# Match the source file: new sheet is the active tab with B16 selected.
$ws.Activate()
$ws.Range("B16").Select()
